$wb = $excel.ActiveWorkbook

$headers = @("S","SC","QS","QSC","BTTS%","Blank%","CS%","FG%","FC%")
$cols = @("J","K","L","M","N","O","P","Q","R")

$homeData = @{
    2 = @(344, 207, 55, 34, 0.5, 0.15, 0.4, 0.55, 0.4)
    3 = @(270, 209, 48, 22, 0.42, 0.16, 0.47, 0.74, 0.21)
    4 = @(348, 183, 65, 19, 0.37, 0.16, 0.47, 0.79, 0.21)
    5 = @(255, 208, 39, 29, 0.37, 0.16, 0.58, 0.63, 0.26)
    6 = @(268, 180, 34, 9, 0.39, 0.17, 0.44, 0.78, 0.22)
    7 = @(307, 155, 44, 29, 0.47, 0.16, 0.37, 0.58, 0.42)
    8 = @(288, 191, 47, 21, 0.6, 0.05, 0.35, 0.65, 0.35)
    9 = @(203, 181, 34, 18, 0.39, 0.28, 0.44, 0.67, 0.22)
    10 = @(352, 168, 39, 23, 0.56, 0.28, 0.28, 0.39, 0.5)
    11 = @(159, 140, 16, 21, 0.62, 0.15, 0.31, 0.46, 0.46)
    12 = @(162, 184, 24, 22, 0.85, 0.15, 0.08, 0.38, 0.54)
    13 = @(229, 237, 30, 28, 0.47, 0.26, 0.37, 0.37, 0.53)
    14 = @(242, 175, 31, 14, 0.35, 0.4, 0.4, 0.35, 0.5)
    15 = @(295, 211, 32, 38, 0.53, 0.32, 0.21, 0.42, 0.53)
    16 = @(199, 251, 19, 36, 0.65, 0.29, 0.12, 0.41, 0.53)
    17 = @(191, 148, 25, 27, 0.77, 0.23, 0.15, 0.23, 0.62)
    18 = @(230, 280, 26, 32, 0.32, 0.53, 0.26, 0.37, 0.53)
    19 = @(257, 275, 28, 37, 0.7, 0.15, 0.15, 0.6, 0.4)
    20 = @(274, 220, 32, 31, 0.55, 0.3, 0.25, 0.55, 0.35)
    21 = @(240, 271, 18, 43, 0.7, 0.2, 0.1, 0.3, 0.7)
}

$awayData = @{
    2 = @(267, 227, 42, 16, 0.37, 0, 0.63, 0.89, 0.11)
    3 = @(266, 282, 42, 25, 0.45, 0.15, 0.45, 0.55, 0.4)
    4 = @(256, 263, 41, 25, 0.68, 0.05, 0.26, 0.58, 0.42)
    5 = @(271, 203, 39, 22, 0.45, 0.35, 0.3, 0.45, 0.45)
    6 = @(226, 293, 31, 38, 0.5, 0.15, 0.35, 0.75, 0.25)
    7 = @(199, 219, 25, 27, 0.33, 0.33, 0.44, 0.44, 0.44)
    8 = @(217, 313, 32, 34, 0.42, 0.37, 0.32, 0.32, 0.58)
    9 = @(230, 235, 36, 25, 0.79, 0.16, 0.16, 0.42, 0.47)
    10 = @(201, 236, 31, 32, 0.61, 0.22, 0.17, 0.5, 0.5)
    11 = @(229, 207, 25, 45, 0.53, 0.32, 0.16, 0.53, 0.47)
    12 = @(195, 215, 30, 24, 0.56, 0.33, 0.22, 0.39, 0.5)
    13 = @(119, 179, 12, 18, 0.54, 0.38, 0.15, 0.31, 0.62)
    14 = @(205, 310, 20, 39, 0.53, 0.37, 0.21, 0.32, 0.58)
    15 = @(186, 250, 20, 41, 0.39, 0.5, 0.22, 0.33, 0.56)
    16 = @(149, 311, 14, 46, 0.5, 0.44, 0.11, 0.39, 0.56)
    17 = @(217, 267, 24, 43, 0.48, 0.43, 0.24, 0.19, 0.67)
    18 = @(241, 282, 28, 46, 0.61, 0.39, 0, 0.28, 0.72)
    19 = @(205, 286, 24, 45, 0.47, 0.47, 0.05, 0.32, 0.68)
    20 = @(126, 195, 13, 36, 0.54, 0.46, 0.08, 0.23, 0.69)
    21 = @(137, 216, 24, 36, 0.69, 0.31, 0.08, 0.23, 0.69)
}

$sheetConfigs = @(
    @{ Name = "home"; Data = $homeData },
    @{ Name = "away"; Data = $awayData }
)

foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.Name)

    # Write new header labels J1:R1
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])1").Value = $headers[$i]
    }

    # Copy the style of the existing header cell (I1) onto the new headers
    $ws.Range("I1").Copy()
    $ws.Range("J1:R1").PasteSpecial(-4122)

    # Write the new stat columns for each data row
    $data = $cfg.Data
    foreach ($row in $data.Keys) {
        $vals = $data[$row]
        for ($i = 0; $i -lt $cols.Length; $i++) {
            $ws.Range("$($cols[$i])$row").Value = $vals[$i]
        }
    }
}

$excel.CutCopyMode = 0
